$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F (reuse the exact header style/format from the other headers)
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# Data for rows 2-13: Timestamp(serial), Seconds, Velocity, Acceleration_SMA, Velocity_Bin, Trening
$data = @(
    @(45686.47265717592, 884.4,  10.4,  0.9649090596607753, "10-15", "Duża Gra"),
    @(45686.47783078704, 1331.4, 14.86, 1.022371675286974,  "10-15", "Duża Gra"),
    @(45686.47819768519, 1363.1, 10.09, 1.327337375708989,  "10-15", "Duża Gra"),
    @(45686.47597662037, 1171.2, 6.79,  2.159049306597029,  "5-10",  "Duża Gra"),
    @(45686.47642800926, 1210.2, 6.19,  2.053060940333777,  "5-10",  "Duża Gra"),
    @(45686.47709699074, 1268,   5.29,  1.802566766738891,  "5-10",  "Duża Gra"),
    @(45686.48738865741, 2157.2, 12.53, 2.920514413288662,  "10-15", "Mała Gra"),
    @(45686.48800555555, 2210.5, 14.14, 3.699334178652081,  "10-15", "Mała Gra"),
    @(45686.49229606482, 2581.2, 11.6,  2.934707062585013,  "10-15", "Mała Gra"),
    @(45686.49086550926, 2457.6, 9.5,   2.75393385546548,   "5-10",  "Mała Gra"),
    @(45686.49257152778, 2605,   9.050000000000001, 2.706500717571804, "5-10", "Mała Gra"),
    @(45686.50051365741, 3291.2, 8.74,  2.74220027242388,   "5-10",  "Mała Gra")
)

# Establish the date-time number format once on the first data cell (this mints both the
# intermediate lower-case numFmt and the final upper-case numFmt/style used for the column).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Apply that same (now-registered) format to the rest of the column in one shot.
$ws.Range("A3:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]

    $row = $row + 1
}

Write-Host "done"
